# Generate Report for Handoff
# Updates the localization-status report for "b.md" (row 3 on every sheet):
# it has now been handed off for translation, so its status moves from
# "Handed back: in sync with en-US" to "Ready for handoff", and the
# zh-cn / de-de sheets get a fresh "Latest Handoff File" + "Latest Handoff
# Datetime" for the new handoff package.

$wb = $excel.ActiveWorkbook

$readyStatus = "Ready for handoff"

# ---- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $readyStatus
$overview.Range("C3").Value = $readyStatus

# ---- Helper data for the per-locale sheets ----------------------------
$locales = @(
    @{ Sheet = "zh-cn"; File = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"; Datetime = "2016-03-08 16:29:22" },
    @{ Sheet = "de-de"; File = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"; Datetime = "2016-03-08 16:29:29" }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Sheet)

    # Status (column B) and Latest Handoff File / Datetime (columns C, D)
    $ws.Range("B3").Value = $readyStatus
    $ws.Range("C3").Value = $locale.File
    $ws.Range("D3").Value = $locale.Datetime

    # Update the hyperlink display text over C3 to match the new file name,
    # keeping the same underlying target address.
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq '$C$3') {
            $h.TextToDisplay = $locale.File
        }
    }
}
